# Applies the "scheduled runner" profit-recalculation update to the
# Golem_Profits workbook: refreshed currentAveragePrice* columns and the
# LeveProfit columns that derive from them, across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Range("H96").Value = 1080
$ws.Range("I96").Value = 325
$ws.Range("K96").Value = 975
$ws.Range("M96").Value = 398
# Row 138
$ws.Range("H138").Value = 1905.5
$ws.Range("I138").Value = 1774.8334
$ws.Range("J138").Value = 2036.1666
$ws.Range("K138").Value = 5324.5002
$ws.Range("L138").Value = 6108.4998
$ws.Range("M138").Value = -184.5002000000004
$ws.Range("N138").Value = -16388.4998
# Row 141
$ws.Range("H141").Value = 768.7143
$ws.Range("I141").Value = 798.3333
$ws.Range("K141").Value = 2394.9999
$ws.Range("M141").Value = 2785.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 597.5
$ws.Range("I20").Value = 597.5
$ws.Range("K20").Value = 597.5
$ws.Range("M20").Value = -350.5
# Row 37
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 1000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -1274
# Row 76
$ws.Range("H76").Value = 19333.666
$ws.Range("J76").Value = 19333.666
$ws.Range("L76").Value = 19333.666
$ws.Range("N76").Value = -19963.666
# Row 79
$ws.Range("H79").Value = 19333.666
$ws.Range("J79").Value = 19333.666
$ws.Range("L79").Value = 19333.666
$ws.Range("N79").Value = -21517.666
# Row 105
$ws.Range("H105").Value = 2111.5557
$ws.Range("I105").Value = 1584
$ws.Range("K105").Value = 1584
$ws.Range("M105").Value = 163

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 18475
$ws.Range("J28").Value = 18475
$ws.Range("L28").Value = 18475
$ws.Range("N28").Value = -18965
# Row 58
$ws.Range("H58").Value = 4811.875
$ws.Range("I58").Value = 4927.857
$ws.Range("K58").Value = 4927.857
$ws.Range("M58").Value = -4724.857
# Row 88
$ws.Range("H88").Value = 29585.25
$ws.Range("J88").Value = 29585.25
$ws.Range("L88").Value = 29585.25
$ws.Range("N88").Value = -30397.25
# Row 91
$ws.Range("H91").Value = 29585.25
$ws.Range("J91").Value = 29585.25
$ws.Range("L91").Value = 29585.25
$ws.Range("N91").Value = -32393.25
# Row 107
$ws.Range("H107").Value = 916.5
$ws.Range("I107").Value = 750
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170
# Row 132
$ws.Range("H132").Value = 1762.5555
$ws.Range("I132").Value = 1277.9412
$ws.Range("K132").Value = 3833.8236
$ws.Range("M132").Value = -1303.8236
# Row 134
$ws.Range("H134").Value = 4237.4
$ws.Range("I134").Value = 2796.75
$ws.Range("K134").Value = 8390.25
$ws.Range("M134").Value = -5855.25
# Row 136
$ws.Range("H136").Value = 4811.875
$ws.Range("I136").Value = 4927.857
$ws.Range("K136").Value = 14783.571
$ws.Range("M136").Value = -12233.571

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
# Row 113
$ws.Range("H113").Value = 281
$ws.Range("I113").Value = 281
$ws.Range("K113").Value = 843
$ws.Range("M113").Value = 1327
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 93
$ws.Range("H93").Value = 62500
$ws.Range("J93").Value = 62500
$ws.Range("L93").Value = 62500
$ws.Range("N93").Value = -66244
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1149.7142
$ws.Range("I132").Value = 1149.7142
$ws.Range("K132").Value = 3449.1426
$ws.Range("M132").Value = -919.1425999999997

$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 5000000
$ws.Range("I14").Value = 5000000
$ws.Range("K14").Value = 5000000
$ws.Range("M14").Value = -4999828
# Row 61
$ws.Range("H61").Value = 1996.7142
$ws.Range("I61").Value = 1996.7142
$ws.Range("K61").Value = 1996.7142
$ws.Range("M61").Value = -1794.7142
# Row 104
$ws.Range("H104").Value = 43444.09
$ws.Range("J104").Value = 43444.09
$ws.Range("L104").Value = 43444.09
$ws.Range("N104").Value = -50432.09
# Row 108
$ws.Range("H108").Value = 45000
$ws.Range("J108").Value = 45000
$ws.Range("L108").Value = 45000
$ws.Range("N108").Value = -52680
# Row 113
$ws.Range("H113").Value = 1996.7142
$ws.Range("I113").Value = 1996.7142
$ws.Range("K113").Value = 1996.7142
$ws.Range("M113").Value = 173.2858000000001
# Row 132
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070

$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 29
$ws.Range("H29").Value = 1005000
$ws.Range("I29").Value = 1005000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1005000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1004710
$ws.Range("N29").ClearContents()
# Row 54
$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29480
# Row 76
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630
# Row 79
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184
# Row 107
$ws.Range("H107").Value = 977.5
$ws.Range("I107").Value = 886.1111
$ws.Range("K107").Value = 2658.3333
$ws.Range("M107").Value = -738.3332999999998
# Row 132
$ws.Range("H132").Value = 1724.6666
$ws.Range("I132").Value = 1724.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5173.9998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2643.9998
$ws.Range("N132").ClearContents()
